$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently looks like:
#   Row 1        -> text column headers ("Lg.", "Threading", "Min.Thread Lg.", ...)
#   Rows 2-162    -> 161 data rows
#
# The target layout is:
#   Row 1         -> numeric column-index header: 0,1,2,...,12
#                    (keeps the bold / centered / bordered style that row 1 already had)
#   Row 2         -> the original text header ("Lg.", "Threading", ...), but with NO
#                    special style (plain, like the rest of the data rows)
#   Rows 3-163    -> the original data rows 2-162, shifted down by exactly one row,
#                    otherwise unchanged.
#
# Implement this by duplicating the existing data block (rows 2-162) one row lower
# (into rows 3-163), which leaves row 1 - and its formatting - completely untouched,
# then overwrite row 1 with the new numeric header and row 2 with the old text
# header (clearing any formatting that row 2 still carries from before).

$ws.Range("A2:M162").Copy()
$ws.Range("A3").PasteSpecial(-4104)   # xlPasteAll

# New numeric header row (row 1), replacing the old text header in place.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10
$ws.Range("L1").Value = 11
$ws.Range("M1").Value = 12

# The old text header now belongs in row 2, with plain (unstyled) formatting.
$ws.Range("A2:M2").ClearFormats()
$ws.Range("A2").Value = "Lg."
$ws.Range("B2").Value = "Threading"
$ws.Range("C2").Value = "Min.Thread Lg."
$ws.Range("D2").Value = "HeadDia."
$ws.Range("E2").Value = "Head Ht."
$ws.Range("F2").Value = "DriveSize"
$ws.Range("G2").Value = "TensileStrength, psi"
$ws.Range("H2").Value = "Specifications Met"
$ws.Range("I2").Value = "Pkg.Qty."
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = "Pkg."
$ws.Range("L2").Value = "thread_size"
$ws.Range("M2").Value = "material_surface"
